{"js": "// Apply the text replacements described by the diff: the header date and\n// the 25 two-digit-by-two-digit multiplication answers in the table.\nconst replacements = [\n  [\"2024-06-23 Sunday\", \"2024-06-24 Monday\"],\n  [\"98\u00d797=9506\", \"93\u00d750=4650\"],\n  [\"61\u00d733=2013\", \"34\u00d712=408\"],\n  [\"63\u00d715=945\", \"77\u00d774=5698\"],\n  [\"67\u00d791=6097\", \"94\u00d766=6204\"],\n  [\"15\u00d733=495\", \"29\u00d777=2233\"],\n  [\"24\u00d775=1800\", \"98\u00d770=6860\"],\n  [\"84\u00d746=3864\", \"38\u00d762=2356\"],\n  [\"11\u00d754=594\", \"59\u00d716=944\"],\n  [\"18\u00d723=414\", \"70\u00d765=4550\"],\n  [\"55\u00d767=3685\", \"41\u00d771=2911\"],\n  [\"91\u00d766=6006\", \"11\u00d746=506\"],\n  [\"87\u00d736=3132\", \"26\u00d733=858\"],\n  [\"81\u00d788=7128\", \"30\u00d724=720\"],\n  [\"49\u00d727=1323\", \"64\u00d717=1088\"],\n  [\"95\u00d752=4940\", \"75\u00d760=4500\"],\n  [\"43\u00d795=4085\", \"99\u00d738=3762\"],\n  [\"21\u00d764=1344\", \"16\u00d775=1200\"],\n  [\"40\u00d795=3800\", \"44\u00d731=1364\"],\n  [\"73\u00d760=4380\", \"37\u00d717=629\"],\n  [\"96\u00d723=2208\", \"24\u00d728=672\"],\n  [\"14\u00d734=476\", \"87\u00d720=1740\"],\n  [\"40\u00d784=3360\", \"47\u00d797=4559\"],\n  [\"98\u00d735=3430\", \"36\u00d727=972\"],\n  [\"33\u00d780=2640\", \"58\u00d771=4118\"],\n  [\"59\u00d737=2183\", \"33\u00d746=1518\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff: the header date and\n# the 25 two-digit-by-two-digit multiplication answers in the table.\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$pairs = @(\n  @(\"2024-06-23 Sunday\", \"2024-06-24 Monday\"),\n  @(\"98\u00d797=9506\", \"93\u00d750=4650\"),\n  @(\"61\u00d733=2013\", \"34\u00d712=408\"),\n  @(\"63\u00d715=945\", \"77\u00d774=5698\"),\n  @(\"67\u00d791=6097\", \"94\u00d766=6204\"),\n  @(\"15\u00d733=495\", \"29\u00d777=2233\"),\n  @(\"24\u00d775=1800\", \"98\u00d770=6860\"),\n  @(\"84\u00d746=3864\", \"38\u00d762=2356\"),\n  @(\"11\u00d754=594\", \"59\u00d716=944\"),\n  @(\"18\u00d723=414\", \"70\u00d765=4550\"),\n  @(\"55\u00d767=3685\", \"41\u00d771=2911\"),\n  @(\"91\u00d766=6006\", \"11\u00d746=506\"),\n  @(\"87\u00d736=3132\", \"26\u00d733=858\"),\n  @(\"81\u00d788=7128\", \"30\u00d724=720\"),\n  @(\"49\u00d727=1323\", \"64\u00d717=1088\"),\n  @(\"95\u00d752=4940\", \"75\u00d760=4500\"),\n  @(\"43\u00d795=4085\", \"99\u00d738=3762\"),\n  @(\"21\u00d764=1344\", \"16\u00d775=1200\"),\n  @(\"40\u00d795=3800\", \"44\u00d731=1364\"),\n  @(\"73\u00d760=4380\", \"37\u00d717=629\"),\n  @(\"96\u00d723=2208\", \"24\u00d728=672\"),\n  @(\"14\u00d734=476\", \"87\u00d720=1740\"),\n  @(\"40\u00d784=3360\", \"47\u00d797=4559\"),\n  @(\"98\u00d735=3430\", \"36\u00d727=972\"),\n  @(\"33\u00d780=2640\", \"58\u00d771=4118\"),\n  @(\"59\u00d737=2183\", \"33\u00d746=1518\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
